# Regenerate save_data to use K (strikeouts) instead of Strike# (count of
# strike pitches) in column G, for each start (rows 2-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 9
    3  = 11
    4  = 7
    5  = 6
    6  = 9
    7  = 7
    8  = 7
    9  = 6
    10 = 8
    11 = 6
    12 = 8
    13 = 9
    14 = 6
    15 = 7
    16 = 7
    17 = 7
    18 = 5
    19 = 8
    20 = 4
    21 = 3
    22 = 6
    23 = 9
    24 = 5
    25 = 6
    26 = 2
    27 = 11
    28 = 6
    29 = 8
    30 = 8
    31 = 10
    32 = 5
    33 = 10
    34 = 5
    35 = 5
    36 = 2
    37 = 3
    38 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
